$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 19:16"

# Row 6: Estados Unidos - updated counts
$ws.Range("B6").Value = 38167
$ws.Range("C6").Value = 13960
$ws.Range("E6").Value = 37593

# Row 8: Alemania - updated counts
$ws.Range("B8").Value = 24806
$ws.Range("C8").Value = 2442
$ws.Range("E8").Value = 24447
$ws.Range("G8").Value = 9
$ws.Range("H8").Value = 93

# Row 12: Suiza - updated counts
$ws.Range("B12").Value = 7367
$ws.Range("C12").Value = 504
$ws.Range("E12").Value = 7138
$ws.Range("G12").Value = 18
$ws.Range("H12").Value = 98

# Rows 20-21: Canada overtakes Dinamarca in ranking (swap + update)
$ws.Range("A20").Value = "Canada"
$ws.Range("B20").Value = 1426
$ws.Range("C20").Value = 98
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = 1392
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 20

$ws.Range("A21").Value = "Dinamarca"
$ws.Range("B21").Value = 1395
$ws.Range("C21").Value = 69
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1381
$ws.Range("F21").Value = 42
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 13

# Row 65: Armenia - updated counts
$ws.Range("B65").Value = 194
$ws.Range("C65").Value = 34
$ws.Range("E65").Value = 192

# Rows 113-114: Nigeria overtakes Guam in ranking (swap + update)
$ws.Range("A113").Value = "Nigeria"
$ws.Range("B113").Value = 30
$ws.Range("C113").Value = 8
$ws.Range("D113").Value = 2
$ws.Range("E113").Value = 28
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A114").Value = "Guam"
$ws.Range("B114").Value = 27
$ws.Range("C114").Value = 12
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 26
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 1
